$p = $ppt.ActivePresentation

# Colour slot order exposed by ThemeColorScheme.Colors(i):
#   1 dk1, 2 lt1, 3 dk2, 4 lt2, 5 accent1, 6 accent2, 7 accent3,
#   8 accent4, 9 accent5, 10 accent6, 11 hlink, 12 folHlink

# "Office Theme" colours -- what the slide-master theme (theme1.xml) becomes.
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

# "Integral" / "Red Violet" colours -- what the notes-master theme (theme2.xml)
# becomes (the deck's previous custom colour set).
$integralColors = @(
    "000000", "FFFFFF", "454551", "D8D9DC",
    "E32D91", "C830CC", "4EA6DC", "4775E7",
    "8971E1", "D54773", "6B9F25", "8C8C8C"
)

# The slide master currently carries the "Integral" theme -- recolour it to the
# standard Office palette.
$masterScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $hex = $officeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $masterScheme.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}

# The notes master currently carries the default "Office Theme" -- recolour it to
# the deck's previous "Integral" palette so that custom look is preserved there.
$notesScheme = $p.NotesMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $hex = $integralColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $notesScheme.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}

# The table on slide 5 switches from the custom "Table_0" style to the built-in
# "Medium Style 2 - Accent 1" table style.
$tableSlide = $p.Slides.Item(5)
for ($i = 1; $i -le $tableSlide.Shapes.Count; $i++) {
    $shp = $tableSlide.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{CDC6D945-6243-4853-B07B-583FA6EE9C67}")
    }
}
